# Akash added new line in AkashFile.docx
#
# Appends a new "PlainText"-styled paragraph (Courier New) after the last
# paragraph, containing two runs:
#   1) "Hi I am adding another line in the AkashFile"
#   2) " and now good to go!"

$d = $word.ActiveDocument

# Insert a brand-new paragraph after the current last paragraph. The new
# paragraph inherits the PlainText style / Courier New run formatting that
# is already in effect at the end of the document.
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newRange = $newPara.Range

# First run.
$newRange.Text = "Hi I am adding another line in the AkashFile"

# Move to the end of what we just typed and append the second run as a
# genuinely separate <w:r> (toggling a character property forces the run
# boundary instead of silently coalescing with the previous run).
$newRange.Collapse(0)
$newRange.InsertAfter(" and now good to go!")
$newRange.Bold = 1
$newRange.Bold = 0
